# LOM3210.xlsx edit script
# Summary of the change (derived from the OOXML diff):
#  - Row 10 (B/C): "Objetivos" description text is replaced by the
#    "519033 - Carlos Yujiro Shigue" docente line.
#  - The two rows that used to hold the "Docentes responsaveis" values
#    (old rows 13 and 14, containing "519033 - Carlos Yujiro Shigue" and
#    "1176388 - Luiz Tadeu Fernandes Eleno") are removed entirely, which
#    shifts every row below them up by two (old row 15 -> new row 13, etc.)
#  - After the shift, several of the remaining descriptive cells are
#    updated in place to their final text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Objetivos" / "Objectives" body text (row 10, cols B & C)
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# 2) Delete the two rows that held the "Docentes responsaveis" names
#    (old rows 13 & 14). Everything below shifts up by two rows.
$ws.Rows("13:14").Delete()

# 3) Fix up the cells whose final text differs from whatever slid into
#    place after the row deletion.
$ws.Range("B13").Value = "01/01/2012"
$ws.Range("C13").Value = "01/01/2012"

$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

$ws.Range("B20").Value = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("C20").Value = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."

$ws.Range("B21").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características da disciplina, não será oferecida recuperação."
